$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 right after
#    the header row, pushing all existing quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# A2 should carry the same style as the other index cells (A3..A9) -
# copy that formatting over before writing the value.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.02

# ------------------------------------------------------------------
# 2) Add a brand-new "2022-Q3" detail sheet, positioned right after
#    "总计" (i.e. right before the current "2022-Q2" sheet). Cloning
#    the existing "2022-Q2" sheet keeps header text/styles identical.
# ------------------------------------------------------------------
$newQ = $wb.Worksheets.Item("2022-Q2")
$newQ.Copy($newQ)
# After Copy(Before:=self), the handle rebinds to the freshly inserted
# clone (now sitting at the original position, original sheet shifted
# one slot later) - rename it in place.
$newQ.Name = "2022-Q3"

$newQ.Range("D2").Value = "0.41"
$newQ.Range("E2").Value = "91.47"
$newQ.Range("F2").Value = "4.66"
$newQ.Range("G2").Value = "0.0191"
